$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.042731434381418
$ws.Range("D2").Value = 1.041576856449796
$ws.Range("E2").Value = 1.040761018512667
$ws.Range("F2").Value = 1.050403512934894
$ws.Range("I2").Value = 1.040744902681959
$ws.Range("J2").Value = 1.047805679147217
$ws.Range("K2").Value = 1.044355978856416
$ws.Range("L2").Value = 1.043542451435393
$ws.Range("M2").Value = 1.053157886400989
$ws.Range("N2").Value = 1.049293683396598
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.044220127057493
$ws.Range("D3").Value = 1.042313192531403
$ws.Range("E3").Value = 1.042046211203316
$ws.Range("F3").Value = 1.052008299216081
$ws.Range("I3").Value = 1.041129731209068
$ws.Range("J3").Value = 1.048938253193758
$ws.Range("K3").Value = 1.044903166044795
$ws.Range("L3").Value = 1.044636885870787
$ws.Range("M3").Value = 1.054573071348753
$ws.Range("N3").Value = 1.050427865828193
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045181788534909
$ws.Range("D4").Value = 1.042788906250701
$ws.Range("E4").Value = 1.042876575032317
$ws.Range("F4").Value = 1.053045414033404
$ws.Range("I4").Value = 1.041376872858322
$ws.Range("J4").Value = 1.049669098038517
$ws.Range("K4").Value = 1.045255822871302
$ws.Range("L4").Value = 1.045343273065088
$ws.Range("M4").Value = 1.055487023359698
$ws.Range("N4").Value = 1.051159748556461
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045585690640847
$ws.Range("D5").Value = 1.042988719446534
$ws.Range("E5").Value = 1.043225368659328
$ws.Range("F5").Value = 1.053481117016812
$ws.Range("I5").Value = 1.041480325948137
$ws.Range("J5").Value = 1.049975870828653
$ws.Range("K5").Value = 1.045403744018809
$ws.Range("L5").Value = 1.045639816389435
$ws.Range("M5").Value = 1.055870833402843
$ws.Range("N5").Value = 1.051466956999146
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.045653485507265
$ws.Range("D6").Value = 1.043022258633653
$ws.Range("E6").Value = 1.043283915732952
$ws.Range("F6").Value = 1.053554256115424
$ws.Range("I6").Value = 1.041497670124055
$ws.Range("J6").Value = 1.050027351683523
$ws.Range("K6").Value = 1.045428560990247
$ws.Range("L6").Value = 1.04568958273989
$ws.Range("M6").Value = 1.05593525266286
$ws.Range("N6").Value = 1.051518510962735
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045187186979378
$ws.Range("D7").Value = 1.042791576857491
$ws.Range("E7").Value = 1.042881236767192
$ws.Range("F7").Value = 1.053051237086547
$ws.Range("I7").Value = 1.041378256950705
$ws.Range("J7").Value = 1.049673199008243
$ws.Range("K7").Value = 1.045257800718175
$ws.Range("L7").Value = 1.045347237142315
$ws.Range("M7").Value = 1.055492153470928
$ws.Range("N7").Value = 1.051163855350033
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.04323488575505
$ws.Range("D8").Value = 1.041825859831666
$ws.Range("E8").Value = 1.041195616115843
$ws.Range("F8").Value = 1.050946130076235
$ws.Range("I8").Value = 1.040875345290038
$ws.Range("J8").Value = 1.048188856845582
$ws.Range("K8").Value = 1.044541196308064
$ws.Range("L8").Value = 1.043912693606118
$ws.Range("M8").Value = 1.053636525360034
$ws.Range("N8").Value = 1.049677405251245
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039781898208861
$ws.Range("D9").Value = 1.040118369307237
$ws.Range("E9").Value = 1.038215558477366
$ws.Range("F9").Value = 1.047226404204474
$ws.Range("I9").Value = 1.039974753721203
$ws.Range("J9").Value = 1.045557623671904
$ws.Range("K9").Value = 1.043267570704804
$ws.Range("L9").Value = 1.041370925716987
$ws.Range("M9").Value = 1.05035278972083
$ws.Range("N9").Value = 1.047042435424654
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.037470807586542
$ws.Range("D10").Value = 1.038976061806128
$ws.Range("E10").Value = 1.036221901388161
$ws.Range("F10").Value = 1.044739150100426
$ws.Range("I10").Value = 1.039364553325774
$ws.Range("J10").Value = 1.043792582388796
$ws.Range("K10").Value = 1.04241106085265
$ws.Range("L10").Value = 1.039666707322353
$ws.Range("M10").Value = 1.048153812572971
$ws.Range("N10").Value = 1.045274887580471
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.036467808887635
$ws.Range("D11").Value = 1.038480463035773
$ws.Range("E11").Value = 1.035356892242461
$ws.Range("F11").Value = 1.043660256992327
$ws.Range("I11").Value = 1.039097975382226
$ws.Range("J11").Value = 1.043025633853671
$ws.Range("K11").Value = 1.042038395393675
$ws.Range("L11").Value = 1.038926381543768
$ws.Range("M11").Value = 1.047199193081839
$ws.Range("N11").Value = 1.044506849890455
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036094897987319
$ws.Range("D12").Value = 1.03829622746782
$ws.Range("E12").Value = 1.035035320161492
$ws.Range("F12").Value = 1.043259212777804
$ws.Range("I12").Value = 1.038998599800528
$ws.Range("J12").Value = 1.042740345935744
$ws.Range("K12").Value = 1.04189969945224
$ws.Range("L12").Value = 1.03865102600142
$ws.Range("M12").Value = 1.046844227666764
$ws.Range("N12").Value = 1.044221156830962
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.036174904751989
$ws.Range("D13").Value = 1.038335753358835
$ws.Range("E13").Value = 1.035104310764176
$ws.Range("F13").Value = 1.043345251654142
$ws.Range("I13").Value = 1.039019932380812
$ws.Range("J13").Value = 1.042801559793033
$ws.Range("K13").Value = 1.041929462521223
$ws.Range("L13").Value = 1.03871010735359
$ws.Range("M13").Value = 1.046920386176586
$ws.Range("N13").Value = 1.04428245761895
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.036436991204656
$ws.Range("D14").Value = 1.038465237099927
$ws.Range("E14").Value = 1.035330316527849
$ws.Range("F14").Value = 1.043627112636887
$ws.Range("I14").Value = 1.0390897682562
$ws.Range("J14").Value = 1.043002060253263
$ws.Range("K14").Value = 1.04202693629379
$ws.Range("L14").Value = 1.038903628066227
$ws.Range("M14").Value = 1.047169859280341
$ws.Range("N14").Value = 1.044483242812831
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.0365984243106
$ws.Range("D15").Value = 1.038544996616652
$ws.Range("E15").Value = 1.035469530248423
$ws.Range("F15").Value = 1.043800737049173
$ws.Range("I15").Value = 1.039132749106363
$ws.Range("J15").Value = 1.043125540775863
$ws.Range("K15").Value = 1.042086957069678
$ws.Range("L15").Value = 1.039022813915666
$ws.Range("M15").Value = 1.047323517621862
$ws.Range("N15").Value = 1.044606898691937
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.037537324409833
$ws.Range("D16").Value = 1.039008932407355
$ws.Range("E16").Value = 1.036279271855716
$ws.Range("F16").Value = 1.044810711851148
$ws.Range("I16").Value = 1.039382195348219
$ws.Range("J16").Value = 1.043843425273437
$ws.Range("K16").Value = 1.042435755499827
$ws.Range("L16").Value = 1.039715789413542
$ws.Range("M16").Value = 1.048217115144493
$ws.Range("N16").Value = 1.045325802667839
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.03812565399873
$ws.Range("D17").Value = 1.039299685363185
$ws.Range("E17").Value = 1.036786729671157
$ws.Range("F17").Value = 1.045443727468949
$ws.Range("I17").Value = 1.039538033532183
$ws.Range("J17").Value = 1.044293014178368
$ws.Range("K17").Value = 1.042654066402794
$ws.Range("L17").Value = 1.040149830511604
$ws.Range("M17").Value = 1.048776983237168
$ws.Range("N17").Value = 1.045776030040584
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03846859725706
$ws.Range("D18").Value = 1.03946918299373
$ws.Range("E18").Value = 1.037082553608892
$ws.Range("F18").Value = 1.045812772576828
$ws.Range("I18").Value = 1.03962870405043
$ws.Range("J18").Value = 1.044554994276534
$ws.Range("K18").Value = 1.042781230807836
$ws.Range("L18").Value = 1.040402769400669
$ws.Range("M18").Value = 1.049103309282761
$ws.Range("N18").Value = 1.046038382180539
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.038585495120395
$ws.Range("D19").Value = 1.039526961482877
$ws.Range("E19").Value = 1.037183393717684
$ws.Range("F19").Value = 1.045938576834046
$ws.Range("I19").Value = 1.039659581883146
$ws.Range("J19").Value = 1.044644279263741
$ws.Range("K19").Value = 1.042824561368427
$ws.Range("L19").Value = 1.040488976228704
$ws.Range("M19").Value = 1.049214538364502
$ws.Range("N19").Value = 1.046127793962669
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.03806255457012
$ws.Range("D20").Value = 1.039268500024372
$ws.Range("E20").Value = 1.036732301641374
$ws.Range("F20").Value = 1.045375829804803
$ws.Range("I20").Value = 1.039521337094576
$ws.Range("J20").Value = 1.044244804231745
$ws.Range("K20").Value = 1.042630661583897
$ws.Range("L20").Value = 1.040103285842699
$ws.Range("M20").Value = 1.048716939120993
$ws.Range("N20").Value = 1.045727751630308
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.036359823087184
$ws.Range("D21").Value = 1.038427111471067
$ws.Range("E21").Value = 1.035263770953503
$ws.Range("F21").Value = 1.043544119838278
$ws.Range("I21").Value = 1.03906921319762
$ws.Range("J21").Value = 1.042943029217061
$ws.Range("K21").Value = 1.041998240207844
$ws.Range("L21").Value = 1.038846651186732
$ws.Range("M21").Value = 1.047096406159082
$ws.Range("N21").Value = 1.044424127945786
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.035287203710453
$ws.Range("D22").Value = 1.037897239076817
$ws.Range("E22").Value = 1.034338886388246
$ws.Range("F22").Value = 1.042390737454307
$ws.Range("I22").Value = 1.03878288013997
$ws.Range("J22").Value = 1.042122181083925
$ws.Range("K22").Value = 1.041599040401966
$ws.Range("L22").Value = 1.038054436519322
$ws.Range("M22").Value = 1.046075325476636
$ws.Range("N22").Value = 1.043602114114149
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035856016624227
$ws.Range("D23").Value = 1.038178216488674
$ws.Range("E23").Value = 1.034829335644613
$ws.Range("F23").Value = 1.043002333228857
$ws.Range("I23").Value = 1.038934867282918
$ws.Range("J23").Value = 1.042557555374883
$ws.Range("K23").Value = 1.041810813500722
$ws.Range("L23").Value = 1.038474607612457
$ws.Range("M23").Value = 1.046616830175427
$ws.Range("N23").Value = 1.044038106686538
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.038091067164052
$ws.Range("D24").Value = 1.039282591626293
$ws.Range("E24").Value = 1.036756895843373
$ws.Range("F24").Value = 1.045406510399752
$ws.Range("I24").Value = 1.039528882198364
$ws.Range("J24").Value = 1.0442665890273
$ws.Range("K24").Value = 1.04264123774779
$ws.Range("L24").Value = 1.040124318085773
$ws.Range("M24").Value = 1.048744071200499
$ws.Range("N24").Value = 1.045749567362771
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040676144818778
$ws.Range("D25").Value = 1.040560490795455
$ws.Range("E25").Value = 1.038987172295075
$ws.Range("F25").Value = 1.048189313491425
$ws.Range("I25").Value = 1.040209297272418
$ws.Range("J25").Value = 1.046239751520525
$ws.Range("K25").Value = 1.043598133983424
$ws.Range("L25").Value = 1.042029718774939
$ws.Range("M25").Value = 1.051203409042699
$ws.Range("N25").Value = 1.047725531973061

Write-Output "Updated 264 cells"